# Apply updated cryptocurrency price/volume figures (and the Aptos/Kaspa row swap)
# to Sheet1, matching the source diff. Values are written as text (NumberFormat "@")
# so Excel does not silently convert numeric-looking strings (e.g. "601.99", "19.31")
# into actual numbers, matching the original inline-string cell typing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.184.40'
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.58%  '
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.623.01'
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.34%  '
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '601.99'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.00%  '
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '150.42'
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +2.79%  '
$ws.Range("E6").Style = "Normal"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.38%  '
$ws.Range("E8").Style = "Normal"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.07%  '
$ws.Range("E9").Style = "Normal"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.74%  '
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.385'
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +6.46%  '
$ws.Range("E11").Style = "Normal"

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.77%  '
$ws.Range("E12").Style = "Normal"

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.88%  '
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.094.26'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.18%  '
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '64.022.01'
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.55%  '
$ws.Range("E15").Style = "Normal"

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +3.64%  '
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.600.99'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.86%  '
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.15'
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +7.53%  '
$ws.Range("E18").Style = "Normal"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +3.88%  '
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '351.17'
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +3.31%  '
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.94'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.13%  '
$ws.Range("E21").Style = "Normal"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.71'
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +2.79%  '
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '66.51'
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.20%  '
$ws.Range("E24").Style = "Normal"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +16.30%  '
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +3.64%  '
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.24'
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +6.86%  '
$ws.Range("E27").Style = "Normal"

$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = 'Aptos'
$ws.Range("B28").Style = "Normal"

$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("C28").Style = "Normal"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.21'
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +4.61%  '
$ws.Range("E28").Style = "Normal"

$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = 'Kaspa'
$ws.Range("B29").Style = "Normal"

$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("C29").Style = "Normal"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.165'
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.50%  '
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.01%  '
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '539.48'
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.60%  '
$ws.Range("E31").Style = "Normal"

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.92%  '
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0₃0854'
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +6.32%  '
$ws.Range("E33").Style = "Normal"

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.78%  '
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.30'
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.54%  '
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '167.85'
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.17%  '
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.409'
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.62%  '
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.98'
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +6.11%  '
$ws.Range("E38").Style = "Normal"

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.06%  '
$ws.Range("E39").Style = "Normal"

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +2.75%  '
$ws.Range("E40").Style = "Normal"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '169.54'
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.99%  '
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.05'
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.28%  '
$ws.Range("E43").Style = "Normal"

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +5.46%  '
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0593'
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +4.10%  '
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '21.63'
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -3.18%  '
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.21%  '
$ws.Range("E47").Style = "Normal"

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.11%  '
$ws.Range("E48").Style = "Normal"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +13.21%  '
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0968'
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.06%  '
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.31'
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +4.10%  '
$ws.Range("E51").Style = "Normal"
